$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the formatting of the last existing header
# cell (AC1) onto the three new header cells, then set their text. ---
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-54): every row gets the same team record. ---
$ws.Range("AD2:AD54").Value = 71
$ws.Range("AE2:AE54").Value = 91
$ws.Range("AF2:AF54").Value = 0

Write-Output $ws.UsedRange.Address()
